# Generate Report for Handback
# Update the status/handback info for the "7fc1593d..." file now that the
# handback has been received and is in sync with en-US.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 is the 7fc1593d...md file
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# zh-cn sheet: row 3 is the 7fc1593d...md file
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$zhcn.Range("K3").Value = "2016-08-27 02:47:57"
$zhcn.Range("P3").Value = ""

# de-de sheet: row 3 is the 7fc1593d...md file
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("K3").Value = "2016-08-27 02:48:09"
$dede.Range("P3").Value = ""

# Error Detail column no longer has long text, so narrow the column to fit.
$zhcn.Columns.Item(16).AutoFit() | Out-Null
$dede.Columns.Item(16).AutoFit() | Out-Null
